$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 2 (row 14): DC vs MI
# Contest 3 (row 15): PBKS vs RCB
# Row 13 is Contest 1 but also had blank score inputs filled in.
# Columns: E=Jaya, H=Justin, K=Ram, N=Sibi, Q=Sundar, T=Upili, W=Vicky (points inputs)

$ws.Range("E13").Value = 60
$ws.Range("H13").Value = 100
$ws.Range("K13").Value = 40
$ws.Range("N13").Value = 0
$ws.Range("Q13").Value = 70
$ws.Range("T13").Value = 50
$ws.Range("W13").Value = 80

$ws.Range("E14").Value = 0
$ws.Range("H14").Value = 70
$ws.Range("K14").Value = 50
$ws.Range("N14").Value = 80
$ws.Range("Q14").Value = 40
$ws.Range("T14").Value = 60
$ws.Range("W14").Value = 100

$ws.Range("E15").Value = 60
$ws.Range("H15").Value = 50
$ws.Range("K15").Value = 40
$ws.Range("N15").Value = 70
$ws.Range("Q15").Value = 80
$ws.Range("T15").Value = 0
$ws.Range("W15").Value = 100

$excel.CalculateFullRebuild()
